$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.118.07"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.526.69"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.13"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.56"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "3.526.46"
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "4.122.00"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.62"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "3.526.90"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "65.048.30"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.09"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  -3.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.68"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "3.667.40"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.71"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -4.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.70"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +8.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.34"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").Value = "3.532.73"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.13"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  +5.47%  "
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.95"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "168.21"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.824"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("E43").Value = "  +4.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.99"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.66"
$ws.Range("E45").Value = "  -5.58%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.91"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "2.416.56"
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.905"
$ws.Range("E51").Value = "  +4.34%  "
